# Updating data for future scenarios
#
# 1) "Import Priorities" sheet (sheet1): add a new "Name variables" header
#    column (A1) and a new explanatory column D with a yellow-highlighted
#    "Don't change this data" header plus two helper notes.
# 2) "Coupling Parameters" sheet (sheet2): bump two forward-looking year /
#    horizon values (End Year 2060 -> 2065, start_tick_fuel_trends 60 -> 70).
# 3) Restore selection/active-sheet state to match the saved workbook.

$wb = $excel.ActiveWorkbook

$wsImport = $wb.Worksheets.Item("Import Priorities")
$wsCoupling = $wb.Worksheets.Item("Coupling Parameters")

# --- Import Priorities: new labelled column of notes -----------------------
# Order matters here: each brand-new distinct string gets appended to the
# shared-string table the first time it is written, so we write them in the
# same left-to-right / top-to-bottom order the authors used.
$wsImport.Range("A1").Value = "Name variables"
$wsImport.Range("D3").Value = "The higher the number the earlier the data is read into the repository."
$wsImport.Range("D2").Value = "Some data needs other data. For example Technologies need fuels. "
$wsImport.Range("D1").Value = "Don’t change this data"

# Highlight the new header cell D1 in yellow.
$wsImport.Range("D1").Interior.Color = 65535

# Widen the new column D so the note text is readable.
$wsImport.Columns.Item(4).ColumnWidth = 38.619791666666664

# --- Coupling Parameters: updated future-scenario figures -------------------
$wsCoupling.Range("B3").Value = 2065
$wsCoupling.Range("B20").Value = 70

# --- Selection / active sheet bookkeeping -----------------------------------
$wsCoupling.Activate()
$wsCoupling.Range("B7").Select()

$wsImport.Activate()
$wsImport.Range("D9").Select()

Write-Host "done"
